$wb = $excel.ActiveWorkbook

$rolesWs = $wb.Worksheets.Item("Roles")
$pagesWs = $wb.Worksheets.Item("Pages")

# Update the "Test Result" row values for HR (C11) and Administrator (F11):
# append the new "-Access when Status Complete" qualifier (introduces a new
# shared string "Yes [ByPassToken]-Access when Status Complete").
$newText = "Yes [ByPassToken]-Access when Status Complete"
$rolesWs.Range("C11").Value = $newText
$rolesWs.Range("F11").Value = $newText

# Widen column C on the Roles sheet (separate sidebar column needs more room).
$rolesWs.Columns.Item(3).ColumnWidth = 25.3

# Roles sheet becomes the active/selected sheet & tab, with a new selection.
$rolesWs.Activate()
$rolesWs.Range("F17").Select()

# Pages sheet keeps its own selection (A8) but is no longer the selected tab.
$pagesWs.Activate()
$pagesWs.Range("A8").Select()

# Re-activate Roles so it is the workbook's active tab/sheet on save.
$rolesWs.Activate()
